# análise descritiva - acréscimo das metas
#
# Insert 5 new columns (meta, meta_avg, meta_std, meta_min, meta_max) right
# before the existing "arrecadado_sucesso" block (column G), shifting all
# later columns (arrecadado_*, apoio_*, contribuicoes_*, menor_ano,
# maior_ano) five positions to the right, and populate the new columns
# with the "meta" statistics for each row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 blank columns at G:K -- this shifts the former G:V block to L:AA
$ws.Range("G1:K1").EntireColumn.Insert()

# The inserted columns copy formatting from the column to their left (F,
# percentage style); restore the proper currency formatting (the format
# used by the former G:O / now L:T "R$ #,##0.00" block) for the new data
# cells.
$ws.Range("G2:K6").NumberFormat = "R$ #,##0.00"

# Header labels for the new columns
$ws.Range("G1").Value = "meta"
$ws.Range("H1").Value = "meta_avg"
$ws.Range("I1").Value = "meta_std"
$ws.Range("J1").Value = "meta_min"
$ws.Range("K1").Value = "meta_max"

# New "meta" statistics data, per row
$ws.Range("G2").Value = 495803.6844660074
$ws.Range("H2").Value = 17096.67877468991
$ws.Range("I2").Value = 15952.79958244942
$ws.Range("J2").Value = 3351.178010772499
$ws.Range("K2").Value = 61717.12812102117

$ws.Range("G3").Value = 2301598.727401814
$ws.Range("H3").Value = 27730.10514941945
$ws.Range("I3").Value = 34657.63167523112
$ws.Range("J3").Value = 46.55761904502517
$ws.Range("K3").Value = 189313.7035611726

$ws.Range("G4").Value = 2221303.500608701
$ws.Range("H4").Value = 15866.45357577644
$ws.Range("I4").Value = 11237.67595303987
$ws.Range("J4").Value = 31.89582864100442
$ws.Range("K4").Value = 80883.37226400203

$ws.Range("G5").Value = 8952422.29225223
$ws.Range("H5").Value = 15542.3998129379
$ws.Range("I5").Value = 13538.7952845963
$ws.Range("J5").Value = 33.25544217501798
$ws.Range("K5").Value = 80687.35013615266

$ws.Range("G6").Value = 1914.395463895744
$ws.Range("H6").Value = 957.1977319478718
$ws.Range("I6").Value = 774.6384084502962
$ws.Range("J6").Value = 409.4456603651128
$ws.Range("K6").Value = 1504.949803530631
